# Feria Lagunitas de Puerto Montt - Zanahoria
# Weekly update: insert two new price-report rows at the top of the data
# block (rows 293-294), shifting the existing rows 293:382 down to 295:384.
#
# This mirrors a new export being prepended to the historical series:
#   - new record dated 2022-06-24 (serial 44736), origin "Chillán"
#   - new record dated 2022-06-24 (serial 44736), origin "Provincia de Llanquihue"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 293, pushing everything
# from 293 downward down by two rows (old 293 -> 295, old 382 -> 384).
$ws.Rows.Item(293).Resize(2).Insert()

# ---- New row 293 ----
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C293").Value = "Los Lagos"
$ws.Range("D293").Value = 44736
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 100114013
$ws.Range("G293").Value = "Zanahoria"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 350
$ws.Range("K293").Value = 10000
$ws.Range("L293").Value = 10000
$ws.Range("M293").Value = 10000
$ws.Range("N293").Value = "$/saco 20 kilos"
$ws.Range("O293").Value = "Chillán"
$ws.Range("P293").Value = 500
$ws.Range("Q293").Value = 20
$ws.Range("R293").Value = "Hortaliza"

# ---- New row 294 ----
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44736
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100114013
$ws.Range("G294").Value = "Zanahoria"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 900
$ws.Range("K294").Value = 7000
$ws.Range("L294").Value = 7000
$ws.Range("M294").Value = 7000
$ws.Range("N294").Value = "$/saco 20 kilos"
$ws.Range("O294").Value = "Provincia de Llanquihue"
$ws.Range("P294").Value = 350
$ws.Range("Q294").Value = 20
$ws.Range("R294").Value = "Hortaliza"

# Date columns use the workbook's date number format (style index 2),
# matching the rest of column D.
$ws.Range("D293:D294").NumberFormat = $ws.Range("D295").NumberFormat
